{"js": "// Replace the multiplication problems' text throughout the document body.\n// Each \"before\" string is unique in the document, so a simple\n// search-and-replace per pair is safe and idempotent.\nconst replacements = [\n  [\"56\u00d779=\", \"27\u00d783=\"],\n  [\"15\u00d718=\", \"79\u00d792=\"],\n  [\"91\u00d742=\", \"33\u00d722=\"],\n  [\"42\u00d744=\", \"59\u00d797=\"],\n  [\"53\u00d799=\", \"58\u00d711=\"],\n  [\"17\u00d762=\", \"40\u00d742=\"],\n  [\"76\u00d792=\", \"35\u00d749=\"],\n  [\"73\u00d798=\", \"69\u00d730=\"],\n  [\"16\u00d779=\", \"96\u00d795=\"],\n  [\"25\u00d744=\", \"83\u00d771=\"],\n  [\"57\u00d778=\", \"45\u00d749=\"],\n  [\"57\u00d715=\", \"96\u00d780=\"],\n  [\"70\u00d768=\", \"92\u00d786=\"],\n  [\"46\u00d740=\", \"16\u00d755=\"],\n  [\"12\u00d791=\", \"38\u00d739=\"],\n  [\"12\u00d743=\", \"20\u00d772=\"],\n  [\"92\u00d773=\", \"48\u00d765=\"],\n  [\"30\u00d745=\", \"53\u00d796=\"],\n  [\"24\u00d789=\", \"13\u00d718=\"],\n  [\"73\u00d752=\", \"52\u00d788=\"],\n  [\"35\u00d759=\", \"94\u00d714=\"],\n  [\"80\u00d711=\", \"66\u00d792=\"],\n  [\"57\u00d763=\", \"32\u00d719=\"],\n  [\"31\u00d755=\", \"67\u00d727=\"],\n  [\"79\u00d729=\", \"14\u00d723=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the multiplication problems' text throughout the document.\n# Each \"before\" string is unique in the document, so a simple\n# Find/Replace pair per entry is safe and idempotent.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"56\u00d779=\", \"27\u00d783=\"),\n    @(\"15\u00d718=\", \"79\u00d792=\"),\n    @(\"91\u00d742=\", \"33\u00d722=\"),\n    @(\"42\u00d744=\", \"59\u00d797=\"),\n    @(\"53\u00d799=\", \"58\u00d711=\"),\n    @(\"17\u00d762=\", \"40\u00d742=\"),\n    @(\"76\u00d792=\", \"35\u00d749=\"),\n    @(\"73\u00d798=\", \"69\u00d730=\"),\n    @(\"16\u00d779=\", \"96\u00d795=\"),\n    @(\"25\u00d744=\", \"83\u00d771=\"),\n    @(\"57\u00d778=\", \"45\u00d749=\"),\n    @(\"57\u00d715=\", \"96\u00d780=\"),\n    @(\"70\u00d768=\", \"92\u00d786=\"),\n    @(\"46\u00d740=\", \"16\u00d755=\"),\n    @(\"12\u00d791=\", \"38\u00d739=\"),\n    @(\"12\u00d743=\", \"20\u00d772=\"),\n    @(\"92\u00d773=\", \"48\u00d765=\"),\n    @(\"30\u00d745=\", \"53\u00d796=\"),\n    @(\"24\u00d789=\", \"13\u00d718=\"),\n    @(\"73\u00d752=\", \"52\u00d788=\"),\n    @(\"35\u00d759=\", \"94\u00d714=\"),\n    @(\"80\u00d711=\", \"66\u00d792=\"),\n    @(\"57\u00d763=\", \"32\u00d719=\"),\n    @(\"31\u00d755=\", \"67\u00d727=\"),\n    @(\"79\u00d729=\", \"14\u00d723=\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $range = $d.Content\n    $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n"}
